# Updated test data for FC test cases
# Adds a "Loading Details Name" / "40V Rail(A)" column (K) to the
# "Add Panels" sheet, duplicating the existing column F header/values,
# and updates the sheet view/row sizing to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Activate()

# Column K gets the same header (row 7) and values (rows 8-9) as column F
# ("Loading Details Name" / "40V Rail(A)"), including its formatting.
$ws.Range("F7:F9").Copy()
$ws.Range("K7:K9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F7:F9").Copy()
$ws.Range("K7:K9").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# New column K should be sized the same as column F
$ws.Columns.Item(11).ColumnWidth = $ws.Columns.Item(6).ColumnWidth()

# Rows 8 and 9 grow to accommodate the wrapped text in the new column
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 28.8

# Scroll the view over and select the newly added column
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$ws.Range("K7:K9").Select()

Write-Host "Added column K (Loading Details Name / 40V Rail(A)) to Add Panels sheet"
